$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the header-row labels: "*_old" -> "*_FV2310" (cols A:J) and
# "*_new" -> "*_FV2404" (cols L:U). Column K ("diff") is left untouched.
$fv2310Headers = @("Segmentname_FV2310", "Segmentgruppe_FV2310", "Segment_FV2310", "Datenelement_FV2310", "Segment ID_FV2310", "Code_FV2310", "Qualifier_FV2310", "Beschreibung_FV2310", "Bedingungsausdruck_FV2310", "Bedingung_FV2310")
for ($i = 0; $i -lt $fv2310Headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 1).Value = $fv2310Headers[$i]
}

$fv2404Headers = @("Segmentname_FV2404", "Segmentgruppe_FV2404", "Segment_FV2404", "Datenelement_FV2404", "Segment ID_FV2404", "Code_FV2404", "Qualifier_FV2404", "Beschreibung_FV2404", "Bedingungsausdruck_FV2404", "Bedingung_FV2404")
for ($i = 0; $i -lt $fv2404Headers.Length; $i++) {
    $ws.Cells.Item(1, $i + 12).Value = $fv2404Headers[$i]
}

# Turn the data range into an Excel Table (ListObject) with autofilter
$lo = $ws.ListObjects.Add(1, $ws.Range("A1:U62"), $false, 1)
$lo.Name = "Table1"

# Freeze the header row (pane split after row 1)
$excel.ActiveWindow.FreezePanes = $false
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
